$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row3 = @(100.88,53.039999999999992,67.14,17.329999999999998,64.400000000000006,71.13,68.84,90.84,78.52,55.25,60.1,49.41,5.07,17.27,48.680000000000007,38.82,-1.1299999999999999,63.289999999999992,50.93,35.690000000000005,3.36,35.299999999999997,22.450000000000003,56.69,20.71,30.09,14.52,23.89,4.3600000000000003,30.800000000000004,17.619999999999997,3.53,2.5100000000000011,5.16,6.4799999999999995,3.05,51.339999999999996,39.069999999999993)
$row4 = @(106.17999999999998,39.840000000000003,58.839999999999996,18.13,81.699999999999989,116.70000000000002,29.290000000000003,104.54,53.22,128.25,47.11,34.53,11.34,61.87,4.2,26.919999999999998,46.790000000000006,67.989999999999995,12.39,51.790000000000006,2.76,26.7,27.35,4.09,37.79,38.79,2.4299999999999997,8.49,0.74,30.800000000000004,3.21,19.05,21.52,0.55000000000000004,26.46,2.0499999999999998,54.239999999999995,32.269999999999996)
$row5 = @(126.17999999999998,35.739999999999995,79.52000000000001,91.240000000000009,91.079999999999984,99.93,7.59,77.47999999999999,64.52,45.5,27.21,14.030000000000001,27.54,7.07,67,27.02,2.17,50.67,45.629999999999995,46.160000000000004,2.56,35,22.450000000000003,30.79,31.39,39.950000000000003,18.12,41.09,16.16,32.300000000000004,15.120000000000001,14.25,12.120000000000001,2.35,5.4,32.049999999999997,97.95999999999998,56.95)
$row6 = @(46.17,72.539999999999992,47.04,29.73,22,72.53,63.04,54.08,33.22,2.81,52.01,14.030000000000001,5.07,55.89,30.279999999999998,9.92,3.5700000000000003,51.69,50.93,16.55,36.36,3.4000000000000004,33.450000000000003,19.490000000000002,29.59,48.45,31.52,38.14,17.04,35,15.120000000000001,5.5500000000000007,12.120000000000001,4.25,26.58,3.05,45.349999999999994,2.0099999999999998)

for ($i = 0; $i -lt $row3.Length; $i++) {
    $col = 3 + $i  # Column C = 3
    $ws.Cells.Item(3, $col).Value = $row3[$i]
    $ws.Cells.Item(4, $col).Value = $row4[$i]
    $ws.Cells.Item(5, $col).Value = $row5[$i]
    $ws.Cells.Item(6, $col).Value = $row6[$i]
}

$ws.Range("B3").Select()

$excel.Calculate()